$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for r=2..15 (Giorno, fcs, TISG_PDR_G, buy_BEE_MWH, sell_lago_MWH, need_to_buy_MW)
$data = @(
    @(45910, 5914.10694019361, 5180.44282366069, 5832, 6521.923268, -1.82253535553832),
    @(45911, 5913.30581487482, 5167.00091212412, 2952, 6521.069439, 117.615189010388),
    @(45912, 5913.14141862252, 4471.63606656278, 2952, 6520.667157, 88.6317418725107),
    @(45913, 1177.4817126383, 2229.38844275892, 2952, 2248.970378, 14.5365461716925),
    @(45914, 1060.75017668403, 2189.53212121692, 2952, 2127.94994, 12.697161855537),
    @(45915, 6061.80298841384, 5015.61062706862, 2952, 6722.443532, 113.510465443949),
    @(45916, 6061.80298841384, 5208.36016751724, 2952, 6722.443532, 121.541696295975),
    @(45917, 6061.80298841384, 5244.4345482863, 2952, 6722.443532, 123.044795494686),
    @(45918, 6061.80298841384, 5214.72781831921, 2952, 6722.443532, 121.80701507939),
    @(45919, 6061.80298841384, 4489.72914577827, 2952, 6722.443532, 91.598737056851),
    @(45920, 1187.6944176276, 2179.15742384782, 2952, 2248.908814, 12.0154925091759),
    @(45921, 1056.51488323291, 2135.46844895674, 2952, 2106.592241, 9.73107528015979),
    @(45922, 5979.99309220113, 4965.66009493116, 2952, 6713.304478, 114.457145030418),
    @(45923, 5979.99309220113, 5214.49650871839, 2952, 6713.304478, 124.825328938219)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
